$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns P and Q (header labels in row 1, values in row 2)
$ws.Range("P1").Value = "Riparian-CanopyCover_score"
$ws.Range("Q1").Value = "Riparian-Disturbance_score"

$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 5
